$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list on Sun Jan 15 20:23:17 UTC 2023 with GitHub Actions
$updates = @{
    "D2" = "300.53"
    "E2" = "-0.86%"
    "E3" = "-1.75%"
    "D4" = "5.117"
    "E4" = "-3.12%"
    "D5" = "0.07365"
    "E5" = "-1.61%"
    "D6" = "2.407"
    "E6" = "61.46%"
    "D7" = "7.940"
    "E7" = "1.07%"
    "D8" = "3.788"
    "E8" = "-0.79%"
    "D9" = "0.9160"
    "E9" = "-0.44%"
    "D10" = "0.1707"
    "E10" = "1.37%"
    "D11" = "0.07578"
    "E11" = "-3.51%"
    "D12" = "0.08091"
    "E12" = "0.73%"
    "D13" = "0.03013"
    "E13" = "-0.20%"
    "D14" = "0.09917"
    "E14" = "0.35%"
    "D15" = "0.001499"
    "E15" = "-0.38%"
    "D16" = "0.006178"
    "E16" = "-2.97%"
    "D17" = "3.468"
    "E17" = "0.20%"
    "D18" = "2.227"
    "E18" = "-0.12%"
    "D19" = "0.3284"
    "E20" = "-0.09%"
    "D21" = "4.642"
    "E21" = "3.37%"
    "D22" = "0.04636"
    "E22" = "0.54%"
    "D23" = "0.1566"
    "E23" = "-3.39%"
    "D24" = "0.001225"
    "E24" = "0.70%"
    "E25" = "0.87%"
    "D26" = "0.0001298"
    "E26" = "-7.10%"
    "E27" = "5.47%"
    "D39" = "0.01716"
    "E39" = "0.10%"
    "D40" = "0.04518"
    "E40" = "0.87%"
    "D41" = "0.007234"
    "E41" = "3.97%"
    "D42" = "0.1343"
    "E42" = "-0.27%"
    "D43" = "0.002227"
    "E43" = "-0.40%"
    "E44" = "-16.49%"
    "E45" = "1.86%"
    "D46" = "0.8085"
    "E46" = "-56.66%"
    "E47" = "-33.21%"
}

foreach ($key in $updates.Keys) {
    $cell = $ws.Range($key)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$key]
    $cell.Style = "Normal"
}
